{"js": "// Remove the \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph along\n// with the blank paragraph before it, the blank paragraph after it, and the\n// following (blank) page-break paragraph -- four consecutive paragraphs in\n// total, right after the \"LOB1036: ... (Requisito fraco)\" requirements line.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nlet markerIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === marker) {\n    markerIndex = i;\n    break;\n  }\n}\n\nif (markerIndex !== -1) {\n  // Delete, from last to first, the paragraph before the marker, the marker\n  // paragraph itself, and the two paragraphs after it.\n  const toDelete = [\n    paragraphs.items[markerIndex + 2], // blank paragraph with pageBreakBefore\n    paragraphs.items[markerIndex + 1], // blank paragraph\n    paragraphs.items[markerIndex],     // \"Ver no Jupiter...\" paragraph\n    paragraphs.items[markerIndex - 1], // blank paragraph before marker\n  ];\n  for (const p of toDelete) {\n    p.delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph along\n# with the blank paragraph before it, the blank paragraph after it, and the\n# following (blank) page-break paragraph -- four consecutive paragraphs in\n# total, right after the \"LOB1036: ... (Requisito fraco)\" requirements line.\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\n$marker = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$idx = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $t = $paras.Item($i).Range.Text.TrimEnd(\"`r\")\n    if ($t -eq $marker) {\n        $idx = $i\n        break\n    }\n}\n\nif ($idx -ge 1) {\n    # Delete from last to first so earlier indices stay valid.\n    $paras.Item($idx + 2).Range.Delete()\n    $paras.Item($idx + 1).Range.Delete()\n    $paras.Item($idx).Range.Delete()\n    $paras.Item($idx - 1).Range.Delete()\n}\n"}
